$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows' D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado) and S (Precio $/Kg) values were shuffled
# between rows 2-13 (the rest of each row's columns stay as-is).
# Mapping: target row <- source row (using the original/before values)
$mapping = @{
    2  = 5
    3  = 8
    4  = 13
    5  = 2
    6  = 3
    7  = 10
    8  = 12
    9  = 4
    10 = 9
    11 = 7
    12 = 6
    13 = 11
}

# Snapshot original values (before any writes) for columns D, M, N, O, P, S
$orig = @{}
foreach ($r in 2..13) {
    $orig[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        S = $ws.Cells.Item($r, 19).Value2
    }
}

foreach ($targetRow in $mapping.Keys) {
    $srcRow = $mapping[$targetRow]
    $src = $orig[$srcRow]

    $ws.Cells.Item($targetRow, 4).Value2  = $src.D
    $ws.Cells.Item($targetRow, 13).Value2 = $src.M
    $ws.Cells.Item($targetRow, 14).Value2 = $src.N
    $ws.Cells.Item($targetRow, 15).Value2 = $src.O
    $ws.Cells.Item($targetRow, 16).Value2 = $src.P
    $ws.Cells.Item($targetRow, 19).Value2 = $src.S
}
